# Apply the "creating function to populate call_history table and fix doc issues"
# edit to the SysCall instances worksheet:
#   1. Fix two mis-typed "locale" values in the existing table (I40, I43 -> SP)
#   2. Add a phone_number / locale lookup table in columns P:Q (rows 36-54),
#      built from the original (pre-fix) phone/locale pairs of the row below
#      each one (i.e. a copy of G38:I55 shifted up one row into P37:Q54,
#      with headers in P36:Q36)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- 1. Fix the two wrong locale values in the main table -------------------
$ws.Range("I40").Value = "SP"
$ws.Range("I43").Value = "SP"

# --- 2. Headers for the new lookup table in P36:Q36 --------------------------
$ws.Range("G37").Copy()
$ws.Range("P36").PasteSpecial(-4122)
$ws.Range("P36").Value = "phone_number"

$ws.Range("I37").Copy()
$ws.Range("Q36").PasteSpecial(-4122)
$ws.Range("Q36").Value = "locale"

# --- 3. Lookup rows P37:Q54 ---------------------------------------------------
# Each row N holds the phone number / locale that originally sat in row N+1's
# G/I columns (before the I40/I43 fix above), i.e. a snapshot of the data.
$phones = @{
    37 = 11985666336
    38 = 11912345678
    39 = 12345678901
    40 = 13945678901
    41 = 14945678901
    42 = 15987654321
    43 = 16976980123
    44 = 17912345678
    45 = 18989456321
    46 = 21934567890
    47 = 22945678901
    48 = 23945678901
    49 = 44934512345
    50 = 45934512345
    51 = 61934512345
    52 = 83912345678
    53 = 99987654321
    54 = 98934561234
}
$locales = @{
    37 = "SP"
    38 = "SP"
    39 = "MG"
    40 = "SP"
    41 = "SP"
    42 = "RJ"
    43 = "SP"
    44 = "SP"
    45 = "SP"
    46 = "RJ"
    47 = "RJ"
    48 = "RJ"
    49 = "PR"
    50 = "PR"
    51 = "DF"
    52 = "PB"
    53 = "RN"
    54 = "AP"
}

$ws.Range("G38").Copy()
$ws.Range("P37:P54").PasteSpecial(-4122)

$ws.Range("I38").Copy()
$ws.Range("Q37:Q54").PasteSpecial(-4122)

foreach ($row in 37..54) {
    $pCell = $ws.Cells.Item($row, 16)
    $pCell.Value = $phones[$row]
    $qCell = $ws.Cells.Item($row, 17)
    $qCell.Value = $locales[$row]
}

# --- 4. Cosmetics: selection + column width for the new column --------------
$ws.Columns.Item(16).ColumnWidth = 13.9
$ws.Range("O52").Select()
